$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (343) down to the new rows
# so the new date cells in column A pick up the same bold/date-format style (s="2").
$ws.Range("A343").Copy($ws.Range("A344:A357"))

$ws.Cells.Item(344, 1).Value = 44418
$ws.Cells.Item(344, 2).Value = 0
$ws.Cells.Item(344, 3).Value = 9
$ws.Cells.Item(344, 4).Value = 56.35566687539136

$ws.Cells.Item(345, 1).Value = 44419
$ws.Cells.Item(345, 2).Value = 0
$ws.Cells.Item(345, 3).Value = 8
$ws.Cells.Item(345, 4).Value = 50.09392611145898

$ws.Cells.Item(346, 1).Value = 44420
$ws.Cells.Item(346, 2).Value = 1
$ws.Cells.Item(346, 3).Value = 5
$ws.Cells.Item(346, 4).Value = 31.30870381966186

$ws.Cells.Item(347, 1).Value = 44421
$ws.Cells.Item(347, 2).Value = 1
$ws.Cells.Item(347, 3).Value = 6
$ws.Cells.Item(347, 4).Value = 37.57044458359425

$ws.Cells.Item(348, 1).Value = 44422
$ws.Cells.Item(348, 2).Value = 0
$ws.Cells.Item(348, 3).Value = 4
$ws.Cells.Item(348, 4).Value = 25.04696305572949

$ws.Cells.Item(349, 1).Value = 44423
$ws.Cells.Item(349, 2).Value = 3
$ws.Cells.Item(349, 3).Value = 6
$ws.Cells.Item(349, 4).Value = 37.57044458359425

$ws.Cells.Item(350, 1).Value = 44424
$ws.Cells.Item(350, 2).Value = 2
$ws.Cells.Item(350, 3).Value = 7
$ws.Cells.Item(350, 4).Value = 43.83218534752661

$ws.Cells.Item(351, 1).Value = 44425
$ws.Cells.Item(351, 2).Value = 0
$ws.Cells.Item(351, 3).Value = 7
$ws.Cells.Item(351, 4).Value = 43.83218534752661

$ws.Cells.Item(352, 1).Value = 44426
$ws.Cells.Item(352, 2).Value = 0
$ws.Cells.Item(352, 3).Value = 7
$ws.Cells.Item(352, 4).Value = 43.83218534752661

$ws.Cells.Item(353, 1).Value = 44427
$ws.Cells.Item(353, 2).Value = 0
$ws.Cells.Item(353, 3).Value = 6
$ws.Cells.Item(353, 4).Value = 37.57044458359425

$ws.Cells.Item(354, 1).Value = 44428
$ws.Cells.Item(354, 2).Value = 4
$ws.Cells.Item(354, 3).Value = 9
$ws.Cells.Item(354, 4).Value = 56.35566687539136

$ws.Cells.Item(355, 1).Value = 44429
$ws.Cells.Item(355, 2).Value = 2
$ws.Cells.Item(355, 3).Value = 11
$ws.Cells.Item(355, 4).Value = 68.8791484032561

$ws.Cells.Item(356, 1).Value = 44430
$ws.Cells.Item(356, 2).Value = 0
$ws.Cells.Item(356, 3).Value = 8
$ws.Cells.Item(356, 4).Value = 50.09392611145898

$ws.Cells.Item(357, 1).Value = 44431
$ws.Cells.Item(357, 2).Value = 4
$ws.Cells.Item(357, 3).Value = 10
$ws.Cells.Item(357, 4).Value = 62.61740763932373

